# Archiving/Unarchiving by Year: the "Group Reviewees" upload-template sheet
# hard-coded the evaluation id "EVAL-001" across row 2 (with a couple of
# cells already reading "EVAL-003"). Update the whole row so every group's
# evaluation id reads the current evaluation, "EVAL-003".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Group Reviewees")

for ($col = 1; $col -le 22; $col++) {
    $ws.Cells.Item(2, $col).Value = "EVAL-003"
}

# Bring the "Group Reviewees" sheet to the front and leave the selection on V2.
$ws.Activate()
$ws.Range("V2").Select()
